# Update TPM-derived ligand/receptor expression & specificity metrics
# for the Adam15-Itgb3 LR-pairs sheet (columns G-T, rows 2-26), per the
# new TPM recomputation described in the commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 48.091872
$ws.Cells.Item(2, 8).Value2 = 144.275616
$ws.Cells.Item(2, 9).Value2 = 0.421093842675958
$ws.Cells.Item(2, 10).Value2 = 0.423782205092405
$ws.Cells.Item(2, 13).Value2 = 6.712486666666666
$ws.Cells.Item(2, 14).Value2 = 20.13746
$ws.Cells.Item(2, 15).Value2 = 0.6330487633990675
$ws.Cells.Item(2, 16).Value2 = 0.6414503882251803
$ws.Cells.Item(2, 17).Value2 = 322.81604957504
$ws.Cells.Item(2, 18).Value2 = 2905.34444617536
$ws.Cells.Item(2, 19).Value2 = 0.2665729363809767
$ws.Cells.Item(2, 20).Value2 = 0.2718352599794462
$ws.Cells.Item(3, 7).Value2 = 48.091872
$ws.Cells.Item(3, 8).Value2 = 144.275616
$ws.Cells.Item(3, 9).Value2 = 0.421093842675958
$ws.Cells.Item(3, 10).Value2 = 0.423782205092405
$ws.Cells.Item(3, 15).Value2 = 0.290741083484562
$ws.Cells.Item(3, 16).Value2 = 0.2945997080427384
$ws.Cells.Item(3, 17).Value2 = 148.260123778944
$ws.Cells.Item(3, 18).Value2 = 1334.341114010496
$ws.Cells.Item(3, 19).Value2 = 0.1224292800682857
$ws.Cells.Item(3, 20).Value2 = 0.1248461138939304
$ws.Cells.Item(4, 7).Value2 = 48.091872
$ws.Cells.Item(4, 8).Value2 = 144.275616
$ws.Cells.Item(4, 9).Value2 = 0.421093842675958
$ws.Cells.Item(4, 10).Value2 = 0.423782205092405
$ws.Cells.Item(4, 13).Value2 = 0.2495096666666667
$ws.Cells.Item(4, 14).Value2 = 0.748529
$ws.Cells.Item(4, 15).Value2 = 0.02353103905946135
$ws.Cells.Item(4, 16).Value2 = 0.02384333563656022
$ws.Cells.Item(4, 17).Value2 = 11.999386952096
$ws.Cells.Item(4, 18).Value2 = 107.994482568864
$ws.Cells.Item(4, 19).Value2 = 0.009908775659706642
$ws.Cells.Item(4, 20).Value2 = 0.01010438135281981
$ws.Cells.Item(5, 7).Value2 = 48.091872
$ws.Cells.Item(5, 8).Value2 = 144.275616
$ws.Cells.Item(5, 9).Value2 = 0.421093842675958
$ws.Cells.Item(5, 10).Value2 = 0.423782205092405
$ws.Cells.Item(5, 13).Value2 = 0.4166465
$ws.Cells.Item(5, 14).Value2 = 0.8332930000000001
$ws.Cells.Item(5, 15).Value2 = 0.03929356804674715
$ws.Cells.Item(5, 16).Value2 = 0.02654337331298611
$ws.Cells.Item(5, 17).Value2 = 20.037310147248
$ws.Cells.Item(5, 18).Value2 = 120.223860883488
$ws.Cells.Item(5, 19).Value2 = 0.016546279561254
$ws.Cells.Item(5, 20).Value2 = 0.01124860927316815
$ws.Cells.Item(6, 7).Value2 = 48.091872
$ws.Cells.Item(6, 8).Value2 = 144.275616
$ws.Cells.Item(6, 9).Value2 = 0.421093842675958
$ws.Cells.Item(6, 10).Value2 = 0.423782205092405
$ws.Cells.Item(6, 11).Value2 = 2
$ws.Cells.Item(6, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(6, 13).Value2 = 0.1419326666666667
$ws.Cells.Item(6, 14).Value2 = 0.425798
$ws.Cells.Item(6, 15).Value2 = 0.01338554601016197
$ws.Cells.Item(6, 16).Value2 = 0.01356319478253491
$ws.Cells.Item(6, 17).Value2 = 6.825807637952001
$ws.Cells.Item(6, 18).Value2 = 61.43226874156801
$ws.Cells.Item(6, 19).Value2 = 0.00563657100573494
$ws.Cells.Item(6, 20).Value2 = 0.005747840593040445
$ws.Cells.Item(7, 9).Value2 = 0.1230362686979479
$ws.Cells.Item(7, 10).Value2 = 0.1238217612582891
$ws.Cells.Item(7, 13).Value2 = 6.712486666666666
$ws.Cells.Item(7, 14).Value2 = 20.13746
$ws.Cells.Item(7, 15).Value2 = 0.6330487633990675
$ws.Cells.Item(7, 16).Value2 = 0.6414503882251803
$ws.Cells.Item(7, 17).Value2 = 94.32121344526222
$ws.Cells.Item(7, 18).Value2 = 848.8909210073598
$ws.Cells.Item(7, 19).Value2 = 0.07788795775247133
$ws.Cells.Item(7, 20).Value2 = 0.07942551682985512
$ws.Cells.Item(8, 9).Value2 = 0.1230362686979479
$ws.Cells.Item(8, 10).Value2 = 0.1238217612582891
$ws.Cells.Item(8, 15).Value2 = 0.290741083484562
$ws.Cells.Item(8, 16).Value2 = 0.2945997080427384
$ws.Cells.Item(8, 19).Value2 = 0.03577169806913908
$ws.Cells.Item(8, 20).Value2 = 0.03647785471602961
$ws.Cells.Item(9, 9).Value2 = 0.1230362686979479
$ws.Cells.Item(9, 10).Value2 = 0.1238217612582891
$ws.Cells.Item(9, 13).Value2 = 0.2495096666666667
$ws.Cells.Item(9, 14).Value2 = 0.748529
$ws.Cells.Item(9, 15).Value2 = 0.02353103905946135
$ws.Cells.Item(9, 16).Value2 = 0.02384333563656022
$ws.Cells.Item(9, 17).Value2 = 3.506011362851555
$ws.Cells.Item(9, 18).Value2 = 31.554102265664
$ws.Cells.Item(9, 19).Value2 = 0.002895171244461795
$ws.Cells.Item(9, 20).Value2 = 0.002952323812791416
$ws.Cells.Item(10, 9).Value2 = 0.1230362686979479
$ws.Cells.Item(10, 10).Value2 = 0.1238217612582891
$ws.Cells.Item(10, 13).Value2 = 0.4166465
$ws.Cells.Item(10, 14).Value2 = 0.8332930000000001
$ws.Cells.Item(10, 15).Value2 = 0.03929356804674715
$ws.Cells.Item(10, 16).Value2 = 0.02654337331298611
$ws.Cells.Item(10, 17).Value2 = 5.854552181514667
$ws.Cells.Item(10, 18).Value2 = 35.127313089088
$ws.Cells.Item(10, 19).Value2 = 0.004834533996300684
$ws.Cells.Item(10, 20).Value2 = 0.003286647233350208
$ws.Cells.Item(11, 9).Value2 = 0.1230362686979479
$ws.Cells.Item(11, 10).Value2 = 0.1238217612582891
$ws.Cells.Item(11, 11).Value2 = 2
$ws.Cells.Item(11, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(11, 13).Value2 = 0.1419326666666667
$ws.Cells.Item(11, 14).Value2 = 0.425798
$ws.Cells.Item(11, 15).Value2 = 0.01338554601016197
$ws.Cells.Item(11, 16).Value2 = 0.01356319478253491
$ws.Cells.Item(11, 17).Value2 = 1.994381815907556
$ws.Cells.Item(11, 18).Value2 = 17.949436343168
$ws.Cells.Item(11, 19).Value2 = 0.001646907635575033
$ws.Cells.Item(11, 20).Value2 = 0.001679418666262709
$ws.Cells.Item(12, 7).Value2 = 21.412221
$ws.Cells.Item(12, 8).Value2 = 64.23666299999999
$ws.Cells.Item(12, 9).Value2 = 0.1874860355013181
$ws.Cells.Item(12, 10).Value2 = 0.1886829905749125
$ws.Cells.Item(12, 13).Value2 = 6.712486666666666
$ws.Cells.Item(12, 14).Value2 = 20.13746
$ws.Cells.Item(12, 15).Value2 = 0.6330487633990675
$ws.Cells.Item(12, 16).Value2 = 0.6414503882251803
$ws.Cells.Item(12, 17).Value2 = 143.72924796622
$ws.Cells.Item(12, 18).Value2 = 1293.56323169598
$ws.Cells.Item(12, 19).Value2 = 0.1186878029287031
$ws.Cells.Item(12, 20).Value2 = 0.1210307775557657
$ws.Cells.Item(13, 7).Value2 = 21.412221
$ws.Cells.Item(13, 8).Value2 = 64.23666299999999
$ws.Cells.Item(13, 9).Value2 = 0.1874860355013181
$ws.Cells.Item(13, 10).Value2 = 0.1886829905749125
$ws.Cells.Item(13, 15).Value2 = 0.290741083484562
$ws.Cells.Item(13, 16).Value2 = 0.2945997080427384
$ws.Cells.Item(13, 17).Value2 = 66.010708334292
$ws.Cells.Item(13, 18).Value2 = 594.096375008628
$ws.Cells.Item(13, 19).Value2 = 0.05450989309987826
$ws.Cells.Item(13, 20).Value2 = 0.05558595393599998
$ws.Cells.Item(14, 7).Value2 = 21.412221
$ws.Cells.Item(14, 8).Value2 = 64.23666299999999
$ws.Cells.Item(14, 9).Value2 = 0.1874860355013181
$ws.Cells.Item(14, 10).Value2 = 0.1886829905749125
$ws.Cells.Item(14, 13).Value2 = 0.2495096666666667
$ws.Cells.Item(14, 14).Value2 = 0.748529
$ws.Cells.Item(14, 15).Value2 = 0.02353103905946135
$ws.Cells.Item(14, 16).Value2 = 0.02384333563656022
$ws.Cells.Item(14, 17).Value2 = 5.342556124302999
$ws.Cells.Item(14, 18).Value2 = 48.08300511872699
$ws.Cells.Item(14, 19).Value2 = 0.004411741224485073
$ws.Cells.Item(14, 20).Value2 = 0.004498831873187568
$ws.Cells.Item(15, 7).Value2 = 21.412221
$ws.Cells.Item(15, 8).Value2 = 64.23666299999999
$ws.Cells.Item(15, 9).Value2 = 0.1874860355013181
$ws.Cells.Item(15, 10).Value2 = 0.1886829905749125
$ws.Cells.Item(15, 13).Value2 = 0.4166465
$ws.Cells.Item(15, 14).Value2 = 0.8332930000000001
$ws.Cells.Item(15, 15).Value2 = 0.03929356804674715
$ws.Cells.Item(15, 16).Value2 = 0.02654337331298611
$ws.Cells.Item(15, 17).Value2 = 8.9213269368765
$ws.Cells.Item(15, 18).Value2 = 53.527961621259
$ws.Cells.Item(15, 19).Value2 = 0.007366995293785893
$ws.Cells.Item(15, 20).Value2 = 0.005008283056640543
$ws.Cells.Item(16, 7).Value2 = 21.412221
$ws.Cells.Item(16, 8).Value2 = 64.23666299999999
$ws.Cells.Item(16, 9).Value2 = 0.1874860355013181
$ws.Cells.Item(16, 10).Value2 = 0.1886829905749125
$ws.Cells.Item(16, 11).Value2 = 2
$ws.Cells.Item(16, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(16, 13).Value2 = 0.1419326666666667
$ws.Cells.Item(16, 14).Value2 = 0.425798
$ws.Cells.Item(16, 15).Value2 = 0.01338554601016197
$ws.Cells.Item(16, 16).Value2 = 0.01356319478253491
$ws.Cells.Item(16, 17).Value2 = 3.039093625786
$ws.Cells.Item(16, 18).Value2 = 27.351842632074
$ws.Cells.Item(16, 19).Value2 = 0.002509602954465753
$ws.Cells.Item(16, 20).Value2 = 0.002559144153318736
$ws.Cells.Item(17, 7).Value2 = 2.1734975
$ws.Cells.Item(17, 8).Value2 = 4.346995
$ws.Cells.Item(17, 9).Value2 = 0.01903120789977957
$ws.Cells.Item(17, 10).Value2 = 0.012768471746644
$ws.Cells.Item(17, 13).Value2 = 6.712486666666666
$ws.Cells.Item(17, 14).Value2 = 20.13746
$ws.Cells.Item(17, 15).Value2 = 0.6330487633990675
$ws.Cells.Item(17, 16).Value2 = 0.6414503882251803
$ws.Cells.Item(17, 17).Value2 = 14.58957298878333
$ws.Cells.Item(17, 18).Value2 = 87.53743793269999
$ws.Cells.Item(17, 19).Value2 = 0.01204768262694602
$ws.Cells.Item(17, 20).Value2 = 0.008190341158927038
$ws.Cells.Item(18, 7).Value2 = 2.1734975
$ws.Cells.Item(18, 8).Value2 = 4.346995
$ws.Cells.Item(18, 9).Value2 = 0.01903120789977957
$ws.Cells.Item(18, 10).Value2 = 0.012768471746644
$ws.Cells.Item(18, 15).Value2 = 0.290741083484562
$ws.Cells.Item(18, 16).Value2 = 0.2945997080427384
$ws.Cells.Item(18, 17).Value2 = 6.700571114870001
$ws.Cells.Item(18, 18).Value2 = 40.20342668922
$ws.Cells.Item(18, 19).Value2 = 0.005533154004801867
$ws.Cells.Item(18, 20).Value2 = 0.003761588048713276
$ws.Cells.Item(19, 7).Value2 = 2.1734975
$ws.Cells.Item(19, 8).Value2 = 4.346995
$ws.Cells.Item(19, 9).Value2 = 0.01903120789977957
$ws.Cells.Item(19, 10).Value2 = 0.012768471746644
$ws.Cells.Item(19, 13).Value2 = 0.2495096666666667
$ws.Cells.Item(19, 14).Value2 = 0.748529
$ws.Cells.Item(19, 15).Value2 = 0.02353103905946135
$ws.Cells.Item(19, 16).Value2 = 0.02384333563656022
$ws.Cells.Item(19, 17).Value2 = 0.5423086367258333
$ws.Cells.Item(19, 18).Value2 = 3.253851820355
$ws.Cells.Item(19, 19).Value2 = 0.0004478240964384426
$ws.Cells.Item(19, 20).Value2 = 0.0003044429574211692
$ws.Cells.Item(20, 7).Value2 = 2.1734975
$ws.Cells.Item(20, 8).Value2 = 4.346995
$ws.Cells.Item(20, 9).Value2 = 0.01903120789977957
$ws.Cells.Item(20, 10).Value2 = 0.012768471746644
$ws.Cells.Item(20, 13).Value2 = 0.4166465
$ws.Cells.Item(20, 14).Value2 = 0.8332930000000001
$ws.Cells.Item(20, 15).Value2 = 0.03929356804674715
$ws.Cells.Item(20, 16).Value2 = 0.02654337331298611
$ws.Cells.Item(20, 17).Value2 = 0.9055801261337501
$ws.Cells.Item(20, 18).Value2 = 3.622320504535
$ws.Cells.Item(20, 19).Value2 = 0.0007478040626217805
$ws.Cells.Item(20, 20).Value2 = 0.0003389183122074875
$ws.Cells.Item(21, 7).Value2 = 2.1734975
$ws.Cells.Item(21, 8).Value2 = 4.346995
$ws.Cells.Item(21, 9).Value2 = 0.01903120789977957
$ws.Cells.Item(21, 10).Value2 = 0.012768471746644
$ws.Cells.Item(21, 11).Value2 = 2
$ws.Cells.Item(21, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(21, 13).Value2 = 0.1419326666666667
$ws.Cells.Item(21, 14).Value2 = 0.425798
$ws.Cells.Item(21, 15).Value2 = 0.01338554601016197
$ws.Cells.Item(21, 16).Value2 = 0.01356319478253491
$ws.Cells.Item(21, 17).Value2 = 0.3084902961683333
$ws.Cells.Item(21, 18).Value2 = 1.85094177701
$ws.Cells.Item(21, 19).Value2 = 0.0002547431089714573
$ws.Cells.Item(21, 20).Value2 = 0.0001731812693750262
$ws.Cells.Item(22, 7).Value2 = 28.477822
$ws.Cells.Item(22, 8).Value2 = 85.433466
$ws.Cells.Item(22, 9).Value2 = 0.2493526452249964
$ws.Cells.Item(22, 10).Value2 = 0.2509445713277496
$ws.Cells.Item(22, 13).Value2 = 6.712486666666666
$ws.Cells.Item(22, 14).Value2 = 20.13746
$ws.Cells.Item(22, 15).Value2 = 0.6330487633990675
$ws.Cells.Item(22, 16).Value2 = 0.6414503882251803
$ws.Cells.Item(22, 17).Value2 = 191.1570004707066
$ws.Cells.Item(22, 18).Value2 = 1720.41300423636
$ws.Cells.Item(22, 19).Value2 = 0.1578523837099703
$ws.Cells.Item(22, 20).Value2 = 0.1609684927011864
$ws.Cells.Item(23, 7).Value2 = 28.477822
$ws.Cells.Item(23, 8).Value2 = 85.433466
$ws.Cells.Item(23, 9).Value2 = 0.2493526452249964
$ws.Cells.Item(23, 10).Value2 = 0.2509445713277496
$ws.Cells.Item(23, 15).Value2 = 0.290741083484562
$ws.Cells.Item(23, 16).Value2 = 0.2945997080427384
$ws.Cells.Item(23, 17).Value2 = 87.79291050834401
$ws.Cells.Item(23, 18).Value2 = 790.1361945750961
$ws.Cells.Item(23, 19).Value2 = 0.07249705824245703
$ws.Cells.Item(23, 20).Value2 = 0.07392819744806516
$ws.Cells.Item(24, 7).Value2 = 28.477822
$ws.Cells.Item(24, 8).Value2 = 85.433466
$ws.Cells.Item(24, 9).Value2 = 0.2493526452249964
$ws.Cells.Item(24, 10).Value2 = 0.2509445713277496
$ws.Cells.Item(24, 13).Value2 = 0.2495096666666667
$ws.Cells.Item(24, 14).Value2 = 0.748529
$ws.Cells.Item(24, 15).Value2 = 0.02353103905946135
$ws.Cells.Item(24, 16).Value2 = 0.02384333563656022
$ws.Cells.Item(24, 17).Value2 = 7.105491874612667
$ws.Cells.Item(24, 18).Value2 = 63.949426871514
$ws.Cells.Item(24, 19).Value2 = 0.005867526834369399
$ws.Cells.Item(24, 20).Value2 = 0.00598335564034026
$ws.Cells.Item(25, 7).Value2 = 28.477822
$ws.Cells.Item(25, 8).Value2 = 85.433466
$ws.Cells.Item(25, 9).Value2 = 0.2493526452249964
$ws.Cells.Item(25, 10).Value2 = 0.2509445713277496
$ws.Cells.Item(25, 13).Value2 = 0.4166465
$ws.Cells.Item(25, 14).Value2 = 0.8332930000000001
$ws.Cells.Item(25, 15).Value2 = 0.03929356804674715
$ws.Cells.Item(25, 16).Value2 = 0.02654337331298611
$ws.Cells.Item(25, 17).Value2 = 11.865184863923
$ws.Cells.Item(25, 18).Value2 = 71.191109183538
$ws.Cells.Item(25, 19).Value2 = 0.009797955132784796
$ws.Cells.Item(25, 20).Value2 = 0.006660915437619727
$ws.Cells.Item(26, 7).Value2 = 28.477822
$ws.Cells.Item(26, 8).Value2 = 85.433466
$ws.Cells.Item(26, 9).Value2 = 0.2493526452249964
$ws.Cells.Item(26, 10).Value2 = 0.2509445713277496
$ws.Cells.Item(26, 11).Value2 = 2
$ws.Cells.Item(26, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(26, 13).Value2 = 0.1419326666666667
$ws.Cells.Item(26, 14).Value2 = 0.425798
$ws.Cells.Item(26, 15).Value2 = 0.01338554601016197
$ws.Cells.Item(26, 16).Value2 = 0.01356319478253491
$ws.Cells.Item(26, 17).Value2 = 4.041933217318667
$ws.Cells.Item(26, 18).Value2 = 36.377398955868
$ws.Cells.Item(26, 19).Value2 = 0.003337721305414783
$ws.Cells.Item(26, 20).Value2 = 0.003403610100537991
